$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the 5 new rows at the correct spots (this shifts existing rows
#    down and preserves their values/formatting automatically, same as a
#    user inserting rows in the live sheet before typing new data into them).
# ---------------------------------------------------------------------------

# After row 4 (Easywood PLA) -> 3 new rows for igus iglidur I150,
# Fillamentum PLA Extrafill, FilaPrint PETG
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# After row 12 (Protopasta Conductive, formerly row 9 + 3 offset) -> 1 new
# row for Proto-pasta HTPLA Matte Fibre
$ws.Rows.Item(13).Insert()

# After row 15 (rigid.ink ABS, formerly row 11 + 4 offset) -> 1 new row for
# rigid.ink ABS Ultra Durable
$ws.Rows.Item(16).Insert()

# ---------------------------------------------------------------------------
# 2) Fix the "Protopast Conductive" -> "Protopasta Conductive" typo (now at
#    row 12 after the inserts above).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Protopasta Conductive"
$ws.Range("B12").Value = "rigid.ink PLA"
$ws.Range("C12").Value = "rigid.ink PLA"

# ---------------------------------------------------------------------------
# 3) Populate the new rows with their formatting (copied from an existing
#    similarly-styled row) and values.
# ---------------------------------------------------------------------------

# Row 5 : igus iglidur I150  (style like row 4/8/9/.. : A=7 B=1 C=0 D=0 E=2)
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A5").Value = "igus iglidur I150"
$ws.Range("B5").Value = "0.15 QUALITY MK3"
$ws.Range("C5").Value = "iglidur I150"
$ws.Range("D5").Value = "MK3 Pretty PETG V2"
$ws.Range("E5").Value = "First one came out with some moderate blistering and some stringing, reduce extrusion multiplier and increase cooling?"

# Row 6 : Fillamentum PLA Extrafill (style A=7 B=0(12) C=0 D=0 E=2)
$ws.Range("A19:E19").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Fillamentum PLA Extrafill"
$ws.Range("B6").Value = "rigid.ink PLA"
$ws.Range("C6").Value = "rigid.ink PLA"
$ws.Range("D6").Value = "Original Prusa i3 MK3 purgebubble"
$ws.Range("E6").Value = ""

# Row 7 : FilaPrint PETG (style A=7 B=1 C=0(12) D=0(12) E=2)
$ws.Range("A4:E4").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A7").Value = "FilaPrint PETG"
$ws.Range("B7").Value = "MK3 Pretty PETG V2"
$ws.Range("C7").Value = "MK3 Pretty PETG V2 filaprint"
$ws.Range("D7").Value = "MK3 Pretty PETG V2"
$ws.Range("E7").Value = ""

# Row 13 : Proto-pasta HTPLA Matte Fibre (style A=7 B=1 C=0 D=0 E=2)
$ws.Range("A4:E4").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Proto-pasta HTPLA Matte Fibre"
$ws.Range("B13").Value = "Proto-pasta HTPLA Matte Fibre"
$ws.Range("C13").Value = "Proto-pasta HTPLA Matte Fibre"
$ws.Range("D13").Value = "Proto-pasta HTPLA Matte Fibre"
$ws.Range("E13").Value = "Consistent, light stringing, first layer 240°C to break nozzle blockage, further layers 210, could potentially go lower."

# Row 16 : rigid.ink ABS Ultra Durable (style A=7 B=0 C=0 D=0 E=2)
$ws.Range("A19:E19").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A16").Value = "rigid.ink ABS Ultra Durable"
$ws.Range("B16").Value = "rigid.ink ABS Ultra Durable"
$ws.Range("C16").Value = "rigid.ink ABS Ultra Durable"
$ws.Range("D16").Value = "MK3 Pretty PETG V2"
$ws.Range("E16").Value = "Need to dry filament before assessing print quality."

# ---------------------------------------------------------------------------
# 4) Column widths: columns A and B are now a uniform ~28.57 (bestFit),
#    replacing the old 22.85546875 / 19.28515625 widths.
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 27.67

# ---------------------------------------------------------------------------
# 5) Selection moves to E8 in the refreshed sheet.
# ---------------------------------------------------------------------------
$ws.Range("E8").Select()
